# Auto-generated from the canonical OOXML diff: updates D/E (and the
# row37/row38 B/C swap) inline-string cell values to the post-edit content.
# A leading apostrophe forces text (quote-prefix) for D-column values that
# would otherwise be auto-coerced to numbers by Excel, matching the source
# string values such as '116.28' or '51.656.60'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.656.60'
$ws.Range('E2').Value = '  +5.33%  '
$ws.Range('D3').Value = '2.759.74'
$ws.Range('E3').Value = '  +5.48%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''116.28'
$ws.Range('E5').Value = '  +4.11%  '
$ws.Range('D6').Value = '''332.47'
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  +2.59%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +6.56%  '
$ws.Range('D10').Value = '''41.77'
$ws.Range('E10').Value = '  +5.45%  '
$ws.Range('E11').Value = '  +6.17%  '
$ws.Range('D12').Value = '''20.18'
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('E14').Value = '  +5.76%  '
$ws.Range('D15').Value = '3.191.06'
$ws.Range('E15').Value = '  +5.47%  '
$ws.Range('D16').Value = '2.731.82'
$ws.Range('E16').Value = '  +4.02%  '
$ws.Range('D17').Value = '''0.887'
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').Value = '51.592.51'
$ws.Range('E18').Value = '  +5.17%  '
$ws.Range('D19').Value = '''3.20'
$ws.Range('E19').Value = '  +6.04%  '
$ws.Range('E20').Value = '  +4.76%  '
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('E22').Value = '  +3.73%  '
$ws.Range('D23').Value = '''278.48'
$ws.Range('E23').Value = '  +3.38%  '
$ws.Range('E24').Value = '  +1.94%  '
$ws.Range('D25').Value = '''2.66'
$ws.Range('E25').Value = '  +5.16%  '
$ws.Range('D26').Value = '''26.83'
$ws.Range('E26').Value = '  +2.93%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +2.44%  '
$ws.Range('D31').Value = '''35.14'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = '''49.94'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('E34').Value = '  +3.03%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('D36').Value = '''19.06'
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''5.00'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '''2.09'
$ws.Range('E38').Value = '  +2.89%  '
$ws.Range('D39').Value = '''3.24'
$ws.Range('E39').Value = '  +3.71%  '
$ws.Range('E40').Value = '  +10.25%  '
$ws.Range('D41').Value = '''127.08'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '''23.12'
$ws.Range('E42').Value = '  +5.07%  '
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('D44').Value = '''2.29'
$ws.Range('E44').Value = '  +7.76%  '
$ws.Range('D45').Value = '''2.44'
$ws.Range('E45').Value = '  +13.01%  '
$ws.Range('D46').Value = '2.091.36'
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('D47').Value = '''3.31'
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('E48').Value = '  +5.00%  '
$ws.Range('D49').Value = '''5.55'
$ws.Range('E49').Value = '  +7.00%  '
$ws.Range('D50').Value = '''9.02'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').Value = '''59.90'
$ws.Range('E51').Value = '  +1.90%  '
